# Relatorios adicionados e avaliacao 2B feita
#
# Marks attendance ("1") in columns G (23/03/2021), H (30/03/2021) and
# I (06/04/2021) for every student row that already had an earlier
# attendance mark, and finally leaves the active selection on B7 (as
# on the sheet after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Presenças ")

# Rows whose students attended on 23/03, 30/03 and 06/04/2021
$fullRows = @(3, 6, 7, 8, 16, 17, 18, 19, 20, 22, 26, 27, 28, 29)

foreach ($r in $fullRows) {
    $ws.Range("G$r").Value = 1
    $ws.Range("H$r").Value = 1
    $ws.Range("I$r").Value = 1
}

# Row 30 only got evaluated for the first two of those dates
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 1

# Leave the selection where the author left it when saving
$ws.Range("B7").Select() | Out-Null
